$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 166; existing rows 166-250 shift down to 167-251
$ws.Rows("166").Insert()

# Populate the newly inserted row 166 with the new record
$ws.Range("A166").Value = 4
$ws.Range("B166").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C166").Value = "Los Lagos"
$ws.Range("D166").Value = 44813
$ws.Range("E166").Value = 10
$ws.Range("F166").Value = "Fruta"
$ws.Range("G166").Value = 100101
$ws.Range("H166").Value = "Berries"
$ws.Range("I166").Value = 100112025
$ws.Range("J166").Value = "Frutilla"
$ws.Range("K166").Value = "Sin especificar"
$ws.Range("L166").Value = "Segunda"
$ws.Range("M166").Value = 100
$ws.Range("N166").Value = 18000
$ws.Range("O166").Value = 18000
$ws.Range("P166").Value = 18000
$ws.Range("Q166").Value = "`$/bandeja 7 kilos"
$ws.Range("R166").Value = "Provincia de Melipilla"
$ws.Range("S166").Value = 2571
$ws.Range("T166").Value = 7
